# Working cosine clustering , still not connecting to filter page
#
# 1) Highlight the rows that the cosine-clustering pass identified as
#    duplicates/near-duplicates of another row (yellow fill across the
#    full row A:V).
# 2) Clean up a few Item_Description (H column) strings on those same rows
#    (and two others) that the clustering pass normalised.
# 3) Refresh a handful of recomputed Unit_Price_USD / TOTAL_ASS_VALUE_USD /
#    Invoice_Unit_Price_FC_USD (T/U/V) numbers that shifted after the
#    clustering re-ran.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) + 2) Rows flagged by the clustering pass: highlight them yellow and
#    tidy up their Item_Description text.
# ---------------------------------------------------------------------

$highlightRows = @(3, 8, 10, 11, 40)
foreach ($r in $highlightRows) {
    $ws.Range("A$r`:V$r").Interior.Color = 65535
}

$ws.Range("H3").Value = "graphitized petroleum coke fc: 985%min s: 005%max ash: 0 50%max vm: 050%max moisture: 050%max size:02-1mm"
$ws.Range("H8").Value = "green(raw) petroleum coke (in bulk)"
$ws.Range("H10").Value = "green(raw) petroleum coke (in bulk)"
$ws.Range("H11").Value = "green(raw) petroleum coke (in bulk)"
$ws.Range("H40").Value = "calcined petroleum coke"

# These two rows get the same description tidy-up but were not part of the
# highlighted cluster.
$ws.Range("H44").Value = "petroleum coke (graphitized)"
$ws.Range("H46").Value = "calcined petroleum coke"

# ---------------------------------------------------------------------
# 3) Recomputed USD columns (T = Unit_Price_USD, U = TOTAL_ASS_VALUE_USD,
#    V = Invoice_Unit_Price_FC_USD) for rows touched by the re-run.
# ---------------------------------------------------------------------

$usdUpdates = @(
    @{ Row = 12; T = 2.2023;      U = 2862.9648;       V = 3.4846 },
    @{ Row = 15; T = 6.3818;      U = 6381553.8694;     V = 0.3487 },
    @{ Row = 16; T = 6.3818;      U = 12763107.4665;    V = 0.3487 },
    @{ Row = 17; T = 6.3818;      U = 15953883.9929;    V = 0.3487 },
    @{ Row = 18; T = 6.3818;      U = 1749471.0113;     V = 0.3487 },
    @{ Row = 19; T = 6.3818;      U = 3190776.7986;     V = 0.3487 },
    @{ Row = 31; T = 2.1837;      U = 3493.902;         V = 3.4846 },
    @{ Row = 34; T = 84.9866;     U = 849841.9051;      V = 0.8075 },
    @{ Row = 35; T = 83.0891;     U = 332356.4061;      V = 0.8075 },
    @{ Row = 36; T = 83.0891;     U = 332356.4061;      V = 0.8075 },
    @{ Row = 37; T = 671.8729;    U = 134373.5613;      V = 8.2282 },
    @{ Row = 38; T = 79928.6032;  U = 1918286.4776;     V = 984.9489 },
    @{ Row = 43; T = 262351.3154; U = 1311756.5536;     V = 2657.6136 }
)

foreach ($u in $usdUpdates) {
    $r = $u.Row
    $ws.Range("T$r").Value = $u.T
    $ws.Range("U$r").Value = $u.U
    $ws.Range("V$r").Value = $u.V
}

Write-Host "cosine clustering edits applied"
